$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each changed coin row.
# NumberFormat is temporarily forced to Text ("@") before writing numeric-looking
# price strings so Excel keeps them as text (matching the source data, which
# stores all Price values as plain strings, e.g. "555.28" not 555.28), then the
# cell style is reset back to Normal so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.134.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.070.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.072.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.03%  "

$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.587.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.248.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.078.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.109"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("E29").Value = "  +2.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "466.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0828"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0403"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.952.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("E43").Value = "  -5.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.72%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("E48").Value = "  +1.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0518"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.29%  "
